$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A96").Value = "2025/12/06 19:00"
$ws.Range("B96").Value = "-"
$ws.Range("C96").Value = "-"
$ws.Range("D96").Value = "-"
$ws.Range("E96").Value = "-"
$ws.Range("F96").Value = "-"
$ws.Range("G96").Value = "-"
